$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$sh = $s.Shapes.Item("Inhaltsplatzhalter 9")
$tf = $sh.TextFrame
$tr = $tf.TextRange

# 1) "Mehrstufiges Vorgehen:" paragraph: remove bold emphasis
$paraVorgehen = $tr.Paragraphs(2)
$paraVorgehen.Font.Bold = 0

# 2) "Erstellung von Regressionsmodellen (je Warengruppe)" paragraph:
#    bold just the word "Regressionsmodellen"
$paraErstellung = $tr.Paragraphs(7)
$paraErstellung.Characters(16, 19).Font.Bold = -1

# 3) "Auswahl des jeweils besten Modells anhand vordefinierter Gütekennzahlen "
#    paragraph: bold "Auswahl" and "besten Modells "
$paraAuswahl = $tr.Paragraphs(8)
$paraAuswahl.Characters(1, 7).Font.Bold = -1
$paraAuswahl.Characters(21, 15).Font.Bold = -1
